$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row 10: name / age / class / etc
$ws.Range("B10").Value = "name"
$ws.Range("C10").Value = "age"
$ws.Range("D10").Value = "class"
$ws.Range("F10").Value = "etc"

# Row 11: Name 1 record
$ws.Range("B11").Value = "Name 1"
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = "2e"
$ws.Range("E11").ClearContents()
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 3

# Row 12: Name 2 record
$ws.Range("B12").Value = "Name 2"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = "3a"
$ws.Range("E12").ClearContents()
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 6

# Update active selection to K10
$ws.Range("K10").Select()
